$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
# Row 6
$ws.Range("H6").Value = 4.33
$ws.Range("O6").Value = 1.29
$ws.Range("P6").Value = 3.75
$ws.Range("Q6").Value = 1.88
$ws.Range("R6").Value = 1.98
$ws.Range("S6").Value = 1.36
$ws.Range("T6").Value = 3
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("X6").Value = 6.5
$ws.Range("AC6").Value = 10
$ws.Range("AD6").Value = 8.5
$ws.Range("AG6").Value = 451
$ws.Range("AK6").Value = 81
$ws.Range("AM6").Value = 51
$ws.Range("AR6").Value = 41
$ws.Range("AT6").Value = 3
$ws.Range("AU6").Value = 9.5
$ws.Range("BB6").Value = 401
# Row 7
$ws.Range("G7").Value = 1.95
$ws.Range("H7").Value = 3.6
$ws.Range("I7").Value = 3.75
$ws.Range("J7").Value = 2.63
$ws.Range("Q7").Value = 2.05
$ws.Range("X7").Value = 9
$ws.Range("Z7").Value = 17
$ws.Range("AC7").Value = 9.5
$ws.Range("AE7").Value = 15
$ws.Range("AH7").Value = 10
$ws.Range("AI7").Value = 19
$ws.Range("AJ7").Value = 13
$ws.Range("AL7").Value = 34
$ws.Range("AO7").Value = 11
$ws.Range("AP7").Value = 21
$ws.Range("AR7").Value = 51
$ws.Range("BB7").Value = 201
# Row 8
$ws.Range("G8").Value = 2.25
$ws.Range("I8").Value = 3.1
$ws.Range("J8").Value = 3
$ws.Range("L8").Value = 3.75
$ws.Range("X8").Value = 11
$ws.Range("Y8").Value = 9.5
$ws.Range("Z8").Value = 21
$ws.Range("AA8").Value = 19
$ws.Range("AB8").Value = 29
$ws.Range("AG8").Value = 201
$ws.Range("AH8").Value = 9.5
$ws.Range("AI8").Value = 15
$ws.Range("AJ8").Value = 11
$ws.Range("AL8").Value = 23
$ws.Range("AO8").Value = 13
$ws.Range("AP8").Value = 23
$ws.Range("AR8").Value = 67
$ws.Range("AW8").Value = 5
$ws.Range("AX8").Value = 17
$ws.Range("BC8").Value = 151
# Row 9
$ws.Range("J9").Value = 3
$ws.Range("L9").Value = 4
$ws.Range("AA9").Value = 19
$ws.Range("AH9").Value = 8.5
$ws.Range("AJ9").Value = 12
$ws.Range("AL9").Value = 29
$ws.Range("AP9").Value = 23
$ws.Range("AV9").Value = 51
# Row 11
$ws.Range("G11").Value = 2.9
$ws.Range("I11").Value = 2.25
$ws.Range("J11").Value = 3.5
$ws.Range("L11").Value = 2.88
$ws.Range("W11").Value = 11
$ws.Range("AA11").Value = 23
$ws.Range("AD11").Value = 7
$ws.Range("AJ11").Value = 9
$ws.Range("AK11").Value = 21
$ws.Range("AL11").Value = 17
$ws.Range("AM11").Value = 23
$ws.Range("AO11").Value = 17
$ws.Range("AX11").Value = 12
# Row 12
$ws.Range("G12").Value = 3.3
$ws.Range("I12").Value = 2.1
$ws.Range("O12").Value = 1.25
$ws.Range("P12").Value = 3.75
$ws.Range("Q12").Value = 1.9
$ws.Range("R12").Value = 1.95
$ws.Range("AQ12").Value = 51
$ws.Range("AS12").Value = 151
$ws.Range("AU12").Value = 7.5
$ws.Range("AW12").Value = 4.33
$ws.Range("AX12").Value = 12
# Row 13
$ws.Range("J13").Value = 4.5
$ws.Range("L13").Value = 2.4
$ws.Range("W13").Value = 12
$ws.Range("AH13").Value = 7.5
# Row 14
$ws.Range("G14").Value = 4.5
$ws.Range("I14").Value = 1.67
$ws.Range("J14").Value = 5
$ws.Range("L14").Value = 2.25
$ws.Range("W14").Value = 15
$ws.Range("AB14").Value = 41
$ws.Range("AX14").Value = 8.5
$ws.Range("AZ14").Value = 26
# Row 15
$ws.Range("O15").Value = 1.36
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 2.15
$ws.Range("R15").Value = 1.67
# Row 16
$ws.Range("I16").Value = 4.2
$ws.Range("J16").Value = 2.5
$ws.Range("L16").Value = 4.75
$ws.Range("Q16").Value = 1.93
$ws.Range("R16").Value = 1.93
$ws.Range("W16").Value = 7
$ws.Range("X16").Value = 8.5
# Row 17
$ws.Range("G17").Value = 2.05
$ws.Range("I17").Value = 3.6
$ws.Range("N17").Value = 8.5
$ws.Range("U17").Value = 1.91
$ws.Range("V17").Value = 1.91
$ws.Range("W17").Value = 7
$ws.Range("X17").Value = 9.5
$ws.Range("Z17").Value = 19
$ws.Range("AG17").Value = 301
$ws.Range("AI17").Value = 17
$ws.Range("AL17").Value = 29
$ws.Range("AV17").Value = 51
$ws.Range("AY17").Value = 29
# Row 18
$ws.Range("G18").Value = 1.65
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = 4.5
$ws.Range("J18").Value = 2.25
$ws.Range("W18").Value = 8
$ws.Range("X18").Value = 8.5
$ws.Range("Z18").Value = 13
$ws.Range("AD18").Value = 7.5
$ws.Range("AL18").Value = 34
$ws.Range("AN18").Value = 3.75
$ws.Range("AO18").Value = 8.5
$ws.Range("AW18").Value = 6.5
# Row 22
$ws.Range("G22").Value = 2.5
$ws.Range("I22").Value = 2.88
$ws.Range("J22").Value = 3.1
$ws.Range("X22").Value = 12
$ws.Range("Z22").Value = 23
$ws.Range("AL22").Value = 23
$ws.Range("AU22").Value = 8
# Row 23
$ws.Range("G23").Value = 1.67
$ws.Range("U23").Value = 1.62
$ws.Range("V23").Value = 2.2
$ws.Range("AE23").Value = 13
$ws.Range("AQ23").Value = 26
$ws.Range("AY23").Value = 26
# Row 24
$ws.Range("G24").Value = 2.63
$ws.Range("I24").Value = 2.75
$ws.Range("J24").Value = 3.4
$ws.Range("M24").Value = 1.07
$ws.Range("N24").Value = 9
$ws.Range("O24").Value = 1.36
$ws.Range("P24").Value = 3
$ws.Range("Q24").Value = 2.2
$ws.Range("R24").Value = 1.65
$ws.Range("Y24").Value = 11
$ws.Range("AG24").Value = 351
